# daily auto push: 2026-01-23 02:29 UTC
# Insert a new data row for 2026/01/23 (Friday) at row 681, shifting all
# subsequent rows down by one (old row 722 becomes row 723).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 681 (pushes existing 681..722 down to 682..723)
$ws.Rows.Item(681).Insert()

# Column A holds dates as plain text (e.g. "2026/01/23"); force text format
# before assigning so Excel does not auto-convert it into a date serial,
# then clear the formatting so the cell keeps the sheet's default style
# (matching the untouched neighboring rows).
$ws.Range("A681").NumberFormat = "@"
$ws.Range("A681").Value = "2026/01/23"
$ws.Range("A681").ClearFormats()

$ws.Range("B681").Value = "金"
$ws.Range("C681").Value = 7
$ws.Range("D681").Value = 18
